$wb = $excel.ActiveWorkbook

# --- "Repayment schedule" sheet: insert a new blank column before column N (14) ---
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Columns.Item(14).Insert()

# --- Make "Repayment schedule" the active sheet/tab and set its selection ---
$ws.Activate()
$ws.Range("L16").Select()

# --- "Transactions" sheet keeps its own selection, just loses the tab-selected flag ---
# (this happens automatically because activating another sheet moves the "active" tab)
